# Mise a jour du mini prototype
# - remove stray spaces around "/" in the "Technologie principale" column
# - make the data rows' font color an explicit black (was theme color)
# - switch the (otherwise unused) column default alignment from General to Left
# - bump a few row heights slightly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Shared-string text fixes ("Python / FastAPI" -> "Python/FastAPI", etc.)
$ws.Range("C2").Value = "Python/FastAPI"
$ws.Range("C4").Value = "Python/ML"
$ws.Range("C5").Value = "Docker/AWS"

# 2) Row heights for the data rows (header row 1 stays as-is)
$ws.Rows.Item(2).RowHeight = 33
$ws.Rows.Item(3).RowHeight = 33
$ws.Rows.Item(4).RowHeight = 19.5
$ws.Rows.Item(5).RowHeight = 19.5

# 3) Column A:D default alignment General -> Left.
#    No visible cell currently uses that default style (every used cell in
#    A1:D5 has its own explicit style), so we touch an empty row that still
#    falls inside columns A:D to flip that shared "column default" format,
#    then clear the scratch row back out so no extra cells/rows remain.
$ws.Range("A6:D6").HorizontalAlignment = -4131
$ws.Range("A6:D6").Clear()

# 4) Data rows font color: was theme-based (automatic/theme 1), now explicit black.
$ws.Range("A2:D5").Font.Color = 0
